$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws.Activate()

$ws.Range("D66").Value = 2
$ws.Rows("67:67").Insert()
